# Applies the "DataSummary" edit described in the commit:
#   "Changes to drive.py to convert to YUV."
# which, on the spreadsheet side, records three new training-run entries
# in the "Models" sheet (rows 32-34) and fixes a duplicated row number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Models")
$ws.Activate()

# ---------------------------------------------------------------------
# Fix duplicated "#" value: row 31 was mistakenly numbered 29 (same as
# row 30); it should be 30.
# ---------------------------------------------------------------------
$ws.Cells.Item(31, 1).Value = 30

# Common values shared by all three new rows (same model/run configuration
# as row 31 just above them).
$modelName = "nvidia1"
$trainingSet = "1,4,6,7,8,9,11,12"
$dataFiltering = "Kept 1/3 of [+/- 0.04]"
$cropping = "[74, 20] [30, 30]"
$dropout = "[0, 0.3, 0.3, 0.5]"
$activation = "relu"

$yellow = 65535

# ---------------------------------------------------------------------
# Row 32 (# 31): steering correction 0.23 test
# ---------------------------------------------------------------------
$r = 32
$ws.Cells.Item($r, 1).Value = 31
$ws.Cells.Item($r, 2).Value = $modelName
$ws.Cells.Item($r, 3).Value = $trainingSet
$ws.Cells.Item($r, 4).Value = 0.23
$ws.Cells.Item($r, 4).Interior.Color = $yellow
$ws.Cells.Item($r, 5).Value = $dataFiltering
$ws.Cells.Item($r, 6).Value = $cropping
$ws.Cells.Item($r, 7).Value = $dropout
$ws.Cells.Item($r, 8).Value = $activation
$ws.Cells.Item($r, 9).Value = 3
$ws.Cells.Item($r, 10).Value = "Stayed along sideline"
$ws.Cells.Item($r, 11).Value = "Confused right off the bat… "
$ws.Cells.Item($r, 12).Value = "More testing of steering angle correction - 0.23 doesn’t seem good.."

# ---------------------------------------------------------------------
# Row 33 (# 32): steering correction 0.27 test
# ---------------------------------------------------------------------
$r = 33
$ws.Cells.Item($r, 1).Value = 32
$ws.Cells.Item($r, 2).Value = $modelName
$ws.Cells.Item($r, 3).Value = $trainingSet
$ws.Cells.Item($r, 4).Value = 0.27
$ws.Cells.Item($r, 4).Interior.Color = $yellow
$ws.Cells.Item($r, 5).Value = $dataFiltering
$ws.Cells.Item($r, 6).Value = $cropping
$ws.Cells.Item($r, 7).Value = $dropout
$ws.Cells.Item($r, 8).Value = $activation
$ws.Cells.Item($r, 9).Value = 3
$ws.Cells.Item($r, 10).Value = "Off road at shadow before bridge"
$ws.Cells.Item($r, 11).Value = "Off road right away"

# ---------------------------------------------------------------------
# Row 34 (# 33): center-only data test
# ---------------------------------------------------------------------
$r = 34
$ws.Cells.Item($r, 1).Value = 33
$ws.Cells.Item($r, 2).Value = $modelName
$ws.Cells.Item($r, 3).Value = $trainingSet
$ws.Cells.Item($r, 4).Value = "Center only"
$ws.Cells.Item($r, 4).Interior.Color = $yellow
$ws.Cells.Item($r, 5).Value = $dataFiltering
$ws.Cells.Item($r, 6).Value = $cropping
$ws.Cells.Item($r, 7).Value = $dropout
$ws.Cells.Item($r, 8).Value = $activation
$ws.Cells.Item($r, 9).Value = 5
$ws.Cells.Item($r, 9).Interior.Color = $yellow
$ws.Cells.Item($r, 10).Value = "Off road at shadow after brdige. Straddled some lines."
$ws.Rows.Item($r).RowHeight = 32

# ---------------------------------------------------------------------
# Update the view so the newly added rows are visible / selected, matching
# the state the workbook was left in after the edit.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D34").Select()
